# Generate Report for Handoff
#
# The localization file "75b38656-db3a-4516-9551-dfd6146543e2.md" has been
# handed off again (status moves from "Handed back: in sync with en-US" to
# "Ready for handoff") and the report rows that track it are refreshed with
# the new handoff timestamps. The zh-cn handback copy on file is now stale,
# so an Error Detail message is recorded for it as well.

$wb = $excel.ActiveWorkbook

$statusReadyForHandoff = "Ready for handoff"
$errorDetailMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/47cc1f08b4cf6ad900d08f0b95810f3658f96d6b/e2e/75b38656-db3a-4516-9551-dfd6146543e2.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/d9a74ac7dcbd8d6ec6ce0f3b7b21699f81a4d7ea/e2e/75b38656-db3a-4516-9551-dfd6146543e2.md."

$dateTimeFormat = "yyyy-mm-dd HH:mm:ss"

# ---- Overview sheet ------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusReadyForHandoff
$wsOverview.Range("F3").Value = $statusReadyForHandoff
$wsOverview.Range("G3").Value = "2016-08-12 08:58:35"
$wsOverview.Range("G3").NumberFormat = $dateTimeFormat

# ---- zh-cn sheet ----------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusReadyForHandoff
$wsZhCn.Range("H3").Value = "2016-08-12 08:58:28"
$wsZhCn.Range("H3").NumberFormat = $dateTimeFormat
$wsZhCn.Range("P3").Value = $errorDetailMessage
$wsZhCn.Range("P1").ColumnWidth = 39.17

# ---- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusReadyForHandoff
$wsDeDe.Range("H3").Value = "2016-08-12 08:58:35"
$wsDeDe.Range("H3").NumberFormat = $dateTimeFormat
$wsDeDe.Range("P3").Value = $errorDetailMessage
$wsDeDe.Range("P1").ColumnWidth = 39.17
